# DeveloperGuide: Rework UML diagrams
#
# The StorageComponentClassDiagram slide modelled an AddressBook app;
# the diagrams are reworked to model a generic Task app instead, so the
# "AddressBook"/"Person" class names become "Task" class names:
#   - AddressBookStorage   -> TaskBookStorage   (<<interface>> box)
#   - XmlAddressBook       -> XmlTaskBook       (class box)
#   - AddressBook          -> TaskBook          (class box, 2nd text line)
#   - XmlAdaptedPerson     -> XmlAdaptedTask    (class box)
#
# NB: the first and third renames keep the original word "Book"/"BookStorage"
# as its own run (only the "Address" prefix becomes "Task"), matching how
# PowerPoint's spell-checker re-flags the freshly split runs (err="1") once
# the replaced word is re-typed. We reproduce that by nudging a (no-op)
# Font property on the trailing substring right after the text swap, which
# forces the engine to materialise it as its own run instead of silently
# merging it back into the preceding one.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- "<<interface>>" / "AddressBookStorage" box -> "TaskBookStorage" ---
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange
$found = $tr.Find("AddressBookStorage")
$found.Text = "TaskBookStorage"
$full = $tr.Text
$splitAt = $full.IndexOf("BookStorage") + 1
$tail = $tr.Characters($splitAt, 11)
$tail.Font.Size = $tail.Font.Size

# --- "XmlAddressBook" / "Storage" box -> "XmlTaskBook" ---
$shp = $s.Shapes.Item(13)
$tr = $shp.TextFrame.TextRange
$found = $tr.Find("XmlAddressBook")
$found.Text = "XmlTaskBook"

# --- "XmlSerializable" / "AddressBook" box -> "TaskBook" ---
$shp = $s.Shapes.Item(20)
$tr = $shp.TextFrame.TextRange
$found = $tr.Find("AddressBook")
$found.Text = "TaskBook"
$full = $tr.Text
$splitAt = $full.IndexOf("Book") + 1
$tail = $tr.Characters($splitAt, 4)
$tail.Font.Size = $tail.Font.Size

# --- "XmlAdaptedPerson" box -> "XmlAdaptedTask" ---
$shp = $s.Shapes.Item(23)
$tr = $shp.TextFrame.TextRange
$found = $tr.Find("XmlAdaptedPerson")
$found.Text = "XmlAdaptedTask"
